# =====================================================================
# Restructure ict-minimal.xlsx from "v1" layout (3 sheets: library_content,
# requirements, scores) to "v2" layout (5 sheets, meta/content split):
#   library_meta, requirements_meta, requirements_content, scores_meta,
#   scores_content
# =====================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: rename the 3 existing sheets to their new v2 names.
# ---------------------------------------------------------------------
$libraryMeta = $wb.Worksheets.Item(1)
$libraryMeta.Name = "library_meta"

$requirementsContent = $wb.Worksheets.Item(2)
$requirementsContent.Name = "requirements_content"

$scoresContent = $wb.Worksheets.Item(3)
$scoresContent.Name = "scores_content"

# ---------------------------------------------------------------------
# Step 2: rewrite library_meta (old "library_content") content.
# Old layout mixed library_*/framework_* rows + tab pointers in one sheet;
# new layout keeps only the "library" key/value block.
# ---------------------------------------------------------------------
$libraryMeta.UsedRange.ClearContents()

$libraryRows = @(
    @("type", "library"),
    @("urn", "urn:intuitem:risk:library:ict-minimal"),
    @("version", "1"),
    @("locale", "en"),
    @("ref_id", "ict-minimal"),
    @("name", "ICT - Minimum standard"),
    @("description", "Minimum standard for improving ICT resilience - Version may 2023"),
    @("copyright", "Creative Commons BY."),
    @("provider", "Swiss FONES"),
    @("packager", "intuitem"),
    @("name[de]", "IKT - Minimalstandard"),
    @("description[de]", "Minimal standard zur Verbesserung der IKT-Resilienz -"),
    @("name[fr]", "Norme minimale TIC"),
    @("description[fr]", "Norme minimale pour améliorer la résilience - version mai 2023"),
    @("name[it]", "Standard minimo TIC"),
    @("description[it]", "Standard minimo per migliorare la resilienza delle TIC")
)

for ($i = 0; $i -lt $libraryRows.Count; $i++) {
    $r = $i + 1
    $libraryMeta.Cells.Item($r, 1).Value = $libraryRows[$i][0]
    $libraryMeta.Cells.Item($r, 2).Value = $libraryRows[$i][1]
}

# ---------------------------------------------------------------------
# Step 3: create the new "requirements_meta" sheet, right after library_meta.
# Holds the framework_* keys (minus the "framework_" prefix) + min/max score
# + a pointer to the scores sheet, that used to live in library_content.
# ---------------------------------------------------------------------
$requirementsMeta = $wb.Worksheets.Add($null, $libraryMeta)
$requirementsMeta.Name = "requirements_meta"

$requirementsMetaRows = @(
    @("type", "framework"),
    @("base_urn", "urn:intuitem:risk:req_node:ict-minimal"),
    @("urn", "urn:intuitem:risk:framework:ict-minimal"),
    @("ref_id", "ict-minimal"),
    @("name", "ICT - Minimum standard"),
    @("description", "Minimum standard for improving ICT resilience - Version may 2023"),
    @("name[de]", "IKT - Minimalstandard"),
    @("description[de]", "Minimal standard zur Verbesserung der IKT-Resilienz"),
    @("name[fr]", "Norme minimale TIC"),
    @("description[fr]", "Norme minimale pour améliorer la résilience - version mai 2023"),
    @("name[it]", "Standard minimo TIC"),
    @("description[it]", "Standard minimo per migliorare la resilienza delle TIC"),
    @("min_score", "0"),
    @("max_score", "4"),
    @("scores_definition", "scores")
)

for ($i = 0; $i -lt $requirementsMetaRows.Count; $i++) {
    $r = $i + 1
    $requirementsMeta.Cells.Item($r, 1).Value = $requirementsMetaRows[$i][0]
    $requirementsMeta.Cells.Item($r, 2).Value = $requirementsMetaRows[$i][1]
}

# ---------------------------------------------------------------------
# Step 4: clean up requirements_content - drop the now-redundant, fully
# empty inlineStr placeholder cells that used to pad every row out to
# column L (a leftover of the v1 generator / cell-merge layout).
# ---------------------------------------------------------------------
$requirementsContentEmptyCells = @(
    "A2","E2","F2","J2","K2","L2","A3","F3","D4","F4","G4","H4","I4","D5","F5",
    "G5","H5","I5","D6","F6","G6","H6","I6","D7","F7","G7","H7","I7","D8","F8",
    "G8","H8","I8","D9","F9","G9","H9","I9","A10","F10","D11","F11","G11","H11","I11",
    "D12","F12","G12","H12","I12","D13","F13","G13","H13","I13","D14","F14","G14","H14","I14",
    "D15","F15","G15","H15","I15","A16","F16","D17","F17","G17","H17","I17","D18","F18","G18",
    "H18","I18","D19","F19","G19","H19","I19","D20","F20","G20","H20","I20","A21","F21","D22",
    "F22","G22","H22","I22","D23","F23","G23","H23","I23","D24","F24","G24","H24","I24","D25",
    "F25","G25","H25","I25","D26","F26","G26","H26","I26","D27","F27","G27","H27","I27","A28",
    "F28","D29","F29","G29","H29","I29","D30","F30","G30","H30","I30","D31","F31","G31","H31",
    "I31","A32","F32","D33","F33","G33","H33","I33","D34","F34","G34","H34","I34","D35","F35",
    "G35","H35","I35","D36","F36","G36","H36","I36","D37","F37","G37","H37","I37","A38","E38",
    "F38","J38","K38","L38","A39","F39","D40","F40","G40","H40","I40","D41","F41","G41","H41",
    "I41","D42","F42","G42","H42","I42","D43","F43","G43","H43","I43","D44","F44","G44","H44",
    "I44","D45","F45","G45","H45","I45","D46","F46","G46","H46","I46","A47","F47","D48","F48",
    "G48","H48","I48","D49","F49","G49","H49","I49","D50","F50","G50","H50","I50","D51","F51",
    "G51","H51","I51","D52","F52","G52","H52","I52","A53","F53","D54","F54","G54","H54","I54",
    "D55","F55","G55","H55","I55","D56","F56","G56","H56","I56","D57","F57","G57","H57","I57",
    "D58","F58","G58","H58","I58","D59","F59","G59","H59","I59","D60","F60","G60","H60","I60",
    "D61","F61","G61","H61","I61","A62","F62","D63","F63","G63","H63","I63","D64","F64","G64",
    "H64","I64","D65","F65","G65","H65","I65","D66","F66","G66","H66","I66","D67","F67","G67",
    "H67","I67","D68","F68","G68","H68","I68","D69","F69","G69","H69","I69","D70","F70","G70",
    "H70","I70","D71","F71","G71","H71","I71","D72","F72","G72","H72","I72","D73","F73","G73",
    "H73","I73","D74","F74","G74","H74","I74","A75","F75","D76","F76","G76","H76","I76","D77",
    "F77","G77","H77","I77","A78","F78","D79","F79","G79","H79","I79","D80","F80","G80","H80",
    "I80","D81","F81","G81","H81","I81","D82","F82","G82","H82","I82","D83","F83","G83","H83",
    "I83","A84","E84","F84","J84","K84","L84","A85","F85","D86","F86","G86","H86","I86","D87",
    "F87","G87","H87","I87","D88","F88","G88","H88","I88","D89","F89","G89","H89","I89","D90",
    "F90","G90","H90","I90","A91","F91","D92","F92","G92","H92","I92","D93","F93","G93","H93",
    "I93","D94","F94","G94","H94","I94","D95","F95","G95","H95","I95","D96","F96","G96","H96",
    "I96","D97","F97","G97","H97","I97","D98","F98","G98","H98","I98","D99","F99","G99","H99",
    "I99","A100","F100","D101","F101","G101","H101","I101","D102","F102","G102","H102","I102","D103","F103",
    "G103","H103","I103","D104","F104","G104","H104","I104","D105","F105","G105","H105","I105","A106","E106",
    "F106","J106","K106","L106","A107","F107","D108","F108","G108","H108","I108","A109","F109","D110","F110",
    "G110","H110","I110","D111","F111","G111","H111","I111","D112","F112","G112","H112","I112","D113","F113",
    "G113","H113","I113","D114","F114","G114","H114","I114","A115","F115","D116","F116","G116","H116","I116",
    "D117","F117","G117","H117","I117","D118","F118","G118","H118","I118","D119","F119","G119","H119","I119",
    "D120","F120","G120","H120","I120","A121","F121","D122","F122","G122","H122","I122","D123","F123","G123",
    "H123","I123","D124","F124","G124","H124","I124","A125","F125","D126","F126","G126","H126","I126","D127",
    "F127","G127","H127","I127","A128","E128","F128","J128","K128","L128","A129","F129","D130","F130","G130",
    "H130","I130","A131","F131","D132","F132","G132","H132","I132","D133","F133","G133","H133","I133","A134",
    "F134","D135","F135","G135","H135","I135","D136","F136","G136","H136","I136","D137","F137","G137","H137",
    "I137"
)
foreach ($addr in $requirementsContentEmptyCells) {
    $requirementsContent.Range($addr).ClearContents()
}

# ---------------------------------------------------------------------
# Step 5: create the new "scores_meta" sheet, right before scores_content.
# ---------------------------------------------------------------------
$scoresContentRef = $wb.Worksheets.Item("scores_content")
$scoresMeta = $wb.Worksheets.Add($scoresContentRef, $null)
$scoresMeta.Name = "scores_meta"

$scoresMeta.Cells.Item(1, 1).Value = "type"
$scoresMeta.Cells.Item(1, 2).Value = "scores"
$scoresMeta.Cells.Item(2, 1).Value = "name"
$scoresMeta.Cells.Item(2, 2).Value = "scores"

# ---------------------------------------------------------------------
# Step 6: clean up scores_content the same way as requirements_content.
# ---------------------------------------------------------------------
$scoresContentEmptyCells = @(
    "C2","E2","G2","I2","C3","E3","G3","I3","C4","E4","G4","I4","C5","E5","G5",
    "I5","C6","E6","G6","I6"
)
foreach ($addr in $scoresContentEmptyCells) {
    $scoresContent.Range($addr).ClearContents()
}

Write-Output "done"
